$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(13.01552393467109,7.105179390761311,0,13.99389279407573,43.87572200638113,48.21658934414275,19.22608362126532,30.91044953836797,9.348973413604547,10.35663680561392,11.9991077916425,15.75078181329377),
    @(12.8493476699279,7.043004225082202,0,14.00689671416528,43.90024702749537,48.28310884820019,19.27366067972739,30.99080627810451,9.351635612720786,10.24259982112736,12.00234146202235,15.73012308262406),
    @(12.74898027932887,7.003759139421257,0,14.01634433714003,43.9239876190918,48.33497924847067,19.30563579043943,31.04502830305453,9.353375771005258,10.17407491453793,12.00582386276471,15.71972584041171),
    @(12.7085473150608,6.987501914223232,0,14.02056302295769,43.93584518889106,48.35888340501894,19.31936031767591,31.06835088446623,9.354111573219825,10.14655703261587,12.00762040632106,15.71606786571892),
    @(12.7018631019317,6.984786565426246,0,14.02128582332941,43.9379459624786,48.36301956135396,19.32168119738311,31.07229763526513,9.354235367469363,10.14201315167831,12.00794154552705,15.71549553667141),
    @(12.74843302922528,7.003540954087469,0,14.01639973793718,43.92413869600127,48.3352904357618,19.30581807296493,31.04533787436671,9.353385586123125,10.17370211237997,12.00584656205761,15.71967415857758),
    @(12.9579101355909,7.0839659021895,0,13.99807334241088,43.88237643559481,48.23723305587787,19.24191465685641,30.93714248861972,9.349869509629292,10.31702542688268,11.99991271578018,15.74318633868688),
    @(13.37974115185466,7.233010281076709,0,13.97371012711399,43.86934372028408,48.13268381110499,19.13853201018193,30.7637622244255,9.34380624546939,10.60855382793796,12.00010635269356,15.80726253230054),
    @(13.69348161179374,7.336991188630661,0,13.96281921141749,43.90166885764146,48.10964074784068,19.07596250363882,30.66009243315521,9.339850718219317,10.82724254982233,12.00739317406243,15.86503546056908),
    @(13.83642811837304,7.383047322023018,0,13.95937566137857,43.92543677070553,48.11087158480235,19.05040757438915,30.6180922365556,9.338157962531916,10.92730580839671,12.01224333101194,15.89357902502025),
    @(13.89054046808256,7.400304461132216,0,13.9582879265756,43.93573553302112,48.11302313371439,19.04114901770496,30.60293080017915,9.337532161682658,10.96524732199063,12.01429918062688,15.90470731812703),
    @(13.87888801661032,7.396596047031406,0,13.95851258756132,43.93345985319275,48.11248479287267,19.04312439378441,30.60616301823177,9.33766626464458,10.95707425876668,12.01384669386097,15.90229652933055),
    @(13.84088062827134,7.384470775150287,0,13.95928184435883,43.92625806629721,48.11101481279162,19.04963747735883,30.61682999628038,9.338106173420401,10.93042647455628,12.01240808330142,15.89448818848088),
    @(13.8175961280606,7.377019698360179,0,13.9597811701391,43.92201568814728,48.1103339166573,19.05368144493291,30.62346063870584,9.338377607219266,10.91410941967985,12.01155539055144,15.88974677890622),
    @(13.68413987923924,7.333955887824331,0,13.96307457307926,43.90029754508213,48.10979613164703,19.0776911293311,30.66294117273251,9.339963478835562,10.82071196759431,12.00710696717053,15.8632151020283),
    @(13.60228985463707,7.307215626885209,0,13.96548122963209,43.88929226704887,48.11246742082949,19.09316541545478,30.68848358044419,9.340963582040599,10.76354020890051,12.00477014971137,15.84751415084876),
    @(13.55523597082742,7.291718369402226,0,13.96700773863386,43.8838157334729,48.11510638902303,19.10233953564725,30.70366048891666,9.341548859185608,10.73071317713231,12.00357068900238,15.83869669415062),
    @(13.53931003380242,7.286451306337352,0,13.96754905158599,43.88210818032615,48.11618920313211,19.10549273608406,30.70888249340375,9.341748752899315,10.71960927462716,12.00318945646975,15.83574806552933),
    @(13.6110008074997,7.310074310041682,0,13.96521032016835,43.89037550506108,48.11206893973075,19.09148981972271,30.68571428196768,9.34085608070855,10.76962062309372,12.00500395195392,15.84916350808764),
    @(13.85204521621245,7.388037269931578,0,13.95905003376014,43.92833821211045,48.11140083574765,19.04771306522369,30.61367666908412,9.337976549771298,10.93825249173829,12.01282470214963,15.89677306548236),
    @(14.00945190797498,7.437920158299548,0,13.95628398649277,43.96071396702425,48.12078800665008,19.02154213488571,30.57092792331427,9.336183202226591,11.04873775081678,12.01921287869886,15.92974782027374),
    @(13.9254692718516,7.411396038223034,0,13.9576453201688,43.94274409340449,48.11487895653278,19.03528671560947,30.59334698923211,9.337132280353002,10.98975558327928,12.01568710343656,15.91198047336269),
    @(13.6070625731524,7.308782285099954,0,13.96533235318868,43.88988312336743,48.11224565663179,19.09224649138014,30.68696474855984,9.340904649944147,10.76687153512928,12.00489780128082,15.84841718151023),
    @(13.26474359345847,7.193642701505028,0,13.97906660715729,43.86550176887806,48.15154074717763,19.16415011383639,30.80650644284534,9.345358295001214,10.52876009470227,11.99879405739892,15.78802990373497)
)

for ($r = 0; $r -lt 24; $r++) {
    for ($c = 0; $c -lt 12; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $data[$r][$c]
    }
}
